$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '
                    <filter>
                        <interfaces xmlns="http://openconfig.net/yang/interfaces">
                        <interface>
                        <name>1/1/1</name>
                        </interface>
                        </interfaces>
                    </filter>
                    '

$ws.Range("F2").Value = '<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:93d20dc6-315e-4c50-9c83-46592205cddd" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <config>
                    <name>1/1/1</name>
                    <type>ethernetCsmacd</type>
                    <mtu>1500</mtu>
                    <description>test</description>
                    <enabled>true</enabled>
                </config>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <config>
                        <auto-negotiate>false</auto-negotiate>
                        <duplex-mode>FULL</duplex-mode>
                        <port-speed>SPEED_100MB</port-speed>
                    </config>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>'

$ws.Range("G2").Value = '  <edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <interfaces xmlns="http://openconfig.net/yang/interfaces">
        <interface>
          <name>1/1/1</name>
          <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
            <config>
              <port-speed>SPEED_1GB</port-speed>
            </config>
          </ethernet>
        </interface>
      </interfaces>
    </config>
  </edit-config>'

$ws.Range("H2").Value = '- Response of edit-config: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:e7175072-9d61-4576-abca-7b810c070b4f" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply> 
 - Response of commit: <?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:4b64ffd3-9aa4-41f0-a967-c06c8f358bad" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <ok/>
</rpc-reply>'

$ws.Range("I2").Value = '<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:4879bc0c-c8bf-4f6a-a3d8-e12475229557" xmlns:nc="urn:ietf:params:xml:ns:netconf:base:1.0" xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
    <data>
        <interfaces xmlns="http://openconfig.net/yang/interfaces">
            <interface>
                <name>1/1/1</name>
                <config>
                    <name>1/1/1</name>
                    <type>ethernetCsmacd</type>
                    <mtu>1500</mtu>
                    <description>test</description>
                    <enabled>true</enabled>
                </config>
                <ethernet xmlns="http://openconfig.net/yang/interfaces/ethernet">
                    <config>
                        <auto-negotiate>false</auto-negotiate>
                        <duplex-mode>FULL</duplex-mode>
                        <port-speed>SPEED_1GB</port-speed>
                    </config>
                </ethernet>
            </interface>
        </interfaces>
    </data>
</rpc-reply>'
